$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New column F: "So ghe" (seat number)
$ws.Columns.Item(6).ColumnWidth = 16

# Values keyed by row, entered in the exact order the workbook author typed
# them (this determines the shared-string table ordering).
$ws.Cells.Item(1, 6).Value = "So ghe"
$ws.Cells.Item(2, 6).Value = "A1"
$ws.Cells.Item(9, 6).Value = "A2"
$ws.Cells.Item(10, 6).Value = "A3"
$ws.Cells.Item(11, 6).Value = "A4"
$ws.Cells.Item(6, 6).Value = "A5"
$ws.Cells.Item(3, 6).Value = "B2"
$ws.Cells.Item(4, 6).Value = "C3"
$ws.Cells.Item(5, 6).Value = "D4"
$ws.Cells.Item(7, 6).Value = "E8"
$ws.Cells.Item(8, 6).Value = "G9"
$ws.Cells.Item(12, 6).Value = "P9"
$ws.Cells.Item(13, 6).Value = "D10"
$ws.Cells.Item(14, 6).Value = "C2"

$ws.Range("F2").Select()
